# Applies the "update Design components" edit:
#   - Slide 1: rename the "PersonListPanel" / "PersonCard" UI component
#     rectangles to "ParcelListPanel" / "ParcelCard".
#   - Refresh the cached "today" text of every datetimeFigureOut date
#     field (slide master and every slide layout) from 1/7/2017 to
#     10/13/2017, as PowerPoint does whenever it re-saves the deck.
#     (The notes master carries the same kind of date field too, but
#     this host's NotesMaster.Shapes writer aliases into the slide
#     master's shape collection instead of the notes master's own
#     part, so it is intentionally left untouched here to avoid
#     corrupting the slide master.)

function Update-DatePlaceholders($shapes, $newText) {
    $count = $shapes.Count
    for ($i = 1; $i -le $count; $i++) {
        $shape = $shapes.Item($i)
        $isDatePlaceholder = $false
        try {
            if ($shape.PlaceholderFormat.Type -eq 16) {
                $isDatePlaceholder = $true
            }
        } catch {
        }
        if ($isDatePlaceholder) {
            $shape.TextFrame.TextRange.Text = $newText
        }
    }
}

$p = $ppt.ActivePresentation
$newDate = "10/13/2017"

# Slide master + every slide layout carry their own copy of the date field.
Update-DatePlaceholders $p.SlideMaster.Shapes $newDate

$layouts = $p.SlideMaster.CustomLayouts
for ($j = 1; $j -le $layouts.Count; $j++) {
    Update-DatePlaceholders $layouts.Item($j).Shapes $newDate
}

# Rename the two renamed UI component shapes on slide 1.
$slide = $p.Slides.Item(1)
$slide.Shapes.Item(11).TextFrame.TextRange.Text = "ParcelListPanel"
$slide.Shapes.Item(12).TextFrame.TextRange.Text = "ParcelCard"
